# Apply the AMFE product-sheet edit: remove the "ETIQUETA" failure-mode
# block (two table rows, originally rows 22-23) from the "AMFE PRODUCTO"
# worksheet. This shifts the following "PÁGINA WEB" block (originally
# rows 24-26, merged B24:B26) up so it becomes rows 22-24 (merged B22:B24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMFE PRODUCTO")
$ws.Activate()

# Delete the two worksheet rows that contain the "ETIQUETA" failure mode
# entries. Excel will automatically shift the rows below upward, fix up
# the merged cell ranges (B24:B26 -> B22:B24), the shared-formula ranges
# and the sheet dimension.
$ws.Rows("22:23").Delete()

# Restore the view/selection state shown in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B22:B24").Select()
